# "Generate Report for Handoff" - refresh the localization status report
# with a new pair of handed-off files (replacing the previous, already
# handed-back pair) and clear out the stale handback bookkeeping columns.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "8c4fc406-04cd-4e0f-adc6-88c7333da1f4"
$oldGuid2 = "8fd3d8bc-dd35-4bc2-bb39-61b7e6ce75db"
$newGuid1 = "021e74fc-61bc-4431-9cca-781eaa3f058e"
$newGuid2 = "ffffebf475d7-e3cf-4ea0-b1e9-3df96f182963"

$newFile1 = "$newGuid1.md"
$newFile2 = "$newGuid2.md"
$newPath1 = "e2e\$newGuid1.md"
$newPath2 = "e2e\$newGuid2.md"

$status = "Ready for handoff"
$hoDateTime = "2016-08-12 17:18:20"

$zhXlf = "$newGuid1.a73c1e716ffe9a8525096c5c6bc0da31c52e85f7.zh-cn.xlf"
$deXlf = "$newGuid1.a73c1e716ffe9a8525096c5c6bc0da31c52e85f7.de-de.xlf"
$zhHandoffDt = "2016-08-12 17:18:13"
$deHandoffDt = "2016-08-12 17:18:20"
$emptyHandback = "0001-01-01 00:00:00"

$oldUrlBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/fd76e91b3810b4529258b89bc7d6ecef14e74180/e2e"

function Set-ColWidth($ws, $colIndex, $targetXmlWidth) {
    $cw = $targetXmlWidth - (5.0/6.0)
    $ws.Columns.Item($colIndex).ColumnWidth = $cw
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value2 = $newFile1
$ov.Range("A3").Value2 = $newFile2
$ov.Range("B2").Value2 = $newPath1
$ov.Range("B3").Value2 = $newPath2

$ov.Range("E2").Value2 = $status
$ov.Range("F2").Value2 = $status
$ov.Range("E3").Value2 = $status
$ov.Range("F3").Value2 = $status

$ov.Range("G2").Value2 = $hoDateTime
$ov.Range("G3").Value2 = $hoDateTime

# rebuild hyperlinks (display text changed); keep original link targets
$ov.Range("A1").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "$oldUrlBase/$oldGuid1.md", "", "", $newPath1) | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "$oldUrlBase/$oldGuid2.md", "", "", $newPath2) | Out-Null

Set-ColWidth $ov 5 17.2159881591797
Set-ColWidth $ov 6 17.2159881591797

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value2 = $newFile1
$zh.Range("A3").Value2 = $newFile2
$zh.Range("C2").Value2 = $status
$zh.Range("C3").Value2 = $status

$zh.Range("G2").Value2 = $zhXlf
$zh.Range("G3").Value2 = $zhXlf
$zh.Range("H2").Value2 = $zhHandoffDt
$zh.Range("H3").Value2 = $zhHandoffDt
$zh.Range("K2").Value2 = $emptyHandback
$zh.Range("K3").Value2 = $emptyHandback

# Content Duplicate (F3) flips False -> True; copy a cell that already
# holds a text "True" shared string so the cell keeps its text type
# instead of Excel auto-coercing the literal into a boolean cell.
$zh.Range("M2").Copy($zh.Range("F3")) | Out-Null

# Latest Target File / Latest Handback File become blank (no handback yet)
$zh.Range("I2").Value2 = ""
$zh.Range("I3").Value2 = ""
$zh.Range("J2").Value2 = ""
$zh.Range("J3").Value2 = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("I3").Style = "Normal"

# rebuild hyperlinks: drop I2/I3 (no longer a handback target), keep A2/A3
$zh.Range("A1").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "$oldUrlBase/$oldGuid1.md", "", "", $newFile1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "$oldUrlBase/$oldGuid2.md", "", "", $newFile2) | Out-Null

Set-ColWidth $zh 3 17.2159881591797
Set-ColWidth $zh 9 18.6506053379604
Set-ColWidth $zh 10 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value2 = $newFile1
$de.Range("A3").Value2 = $newFile2

$de.Range("G2").Value2 = $deXlf
$de.Range("G3").Value2 = $deXlf
$de.Range("H2").Value2 = $deHandoffDt
$de.Range("H3").Value2 = $deHandoffDt
$de.Range("K2").Value2 = $emptyHandback
$de.Range("K3").Value2 = $emptyHandback

$de.Range("M2").Copy($de.Range("F3")) | Out-Null

$de.Range("I2").Value2 = ""
$de.Range("I3").Value2 = ""
$de.Range("J2").Value2 = ""
$de.Range("J3").Value2 = ""
$de.Range("I2").Style = "Normal"
$de.Range("I3").Style = "Normal"

$de.Range("A1").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "$oldUrlBase/$oldGuid1.md", "", "", $newFile1) | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "$oldUrlBase/$oldGuid2.md", "", "", $newFile2) | Out-Null

Set-ColWidth $de 3 17.2159881591797
Set-ColWidth $de 9 18.6506053379604
Set-ColWidth $de 10 21.7054770333426
